$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the requisite headers: PREREQ_COURSES -> Prerequisites, COREQ_COURSES -> CoRequisites
$ws.Range("E1").Value = "Prerequisites"
$ws.Range("F1").Value = "CoRequisites"

# Collapse the header-row selection down to just the header row
[void]$ws.Range("A1:I1").Select()
